# Update gh-pages output data values in sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1488
$ws1.Range("F5").Value = 34
$ws1.Range("F7").Value = 663
$ws1.Range("F14").Value = 156
$ws1.Range("F18").Value = 107
$ws1.Range("F19").Value = 4951
$ws1.Range("F21").Value = 831
$ws1.Range("F23").Value = 2240
$ws1.Range("F26").Value = 2090

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1488
$ws4.Range("F5").Value = 34
$ws4.Range("F7").Value = 663
$ws4.Range("F14").Value = 156
$ws4.Range("F18").Value = 107
$ws4.Range("F19").Value = 4951
$ws4.Range("F23").Value = 831
$ws4.Range("F25").Value = 2240
$ws4.Range("F28").Value = 2090
